$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 3) matching the diff
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Dave"
$ws.Range("C3").Value = "Robinson "
$ws.Range("D3").Value = "robinsondave_876@yahoo.com"
$ws.Range("E3").Value = 96.8
$ws.Range("F3").Value = 24
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = $false
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = "Not at Risk"
